$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 212.4076363333333
$ws.Range("H2").Value = 637.222909
$ws.Range("I2").Value = 0.5080632835800084
$ws.Range("J2").Value = 0.5080632835800084
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 169.629438
$ws.Range("N2").Value = 508.888314
$ws.Range("O2").Value = 0.7428377317484701
$ws.Range("P2").Value = 0.7428377317484702
$ws.Range("Q2").Value = 36030.58797813171
$ws.Range("R2").Value = 324275.2918031854
$ws.Range("S2").Value = 0.3774085771592532
$ws.Range("T2").Value = 0.3774085771592532

$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 212.4076363333333
$ws.Range("H3").Value = 637.222909
$ws.Range("I3").Value = 0.5080632835800084
$ws.Range("J3").Value = 0.5080632835800084
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.9848756666666668
$ws.Range("N3").Value = 2.954627
$ws.Range("O3").Value = 0.004312947180081616
$ws.Range("P3").Value = 0.004312947180081616
$ws.Range("Q3").Value = 209.1951124388826
$ws.Range("R3").Value = 1882.756011949943
$ws.Range("S3").Value = 0.002191250106219404
$ws.Range("T3").Value = 0.002191250106219404

$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 212.4076363333333
$ws.Range("H4").Value = 637.222909
$ws.Range("I4").Value = 0.5080632835800084
$ws.Range("J4").Value = 0.5080632835800084
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 54.620752
$ws.Range("N4").Value = 163.862256
$ws.Range("O4").Value = 0.2391940691454494
$ws.Range("P4").Value = 0.2391940691454494
$ws.Range("Q4").Value = 11601.86482706919
$ws.Range("R4").Value = 104416.7834436227
$ws.Range("S4").Value = 0.1215257241829006
$ws.Range("T4").Value = 0.1215257241829006

$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 212.4076363333333
$ws.Range("H5").Value = 637.222909
$ws.Range("I5").Value = 0.5080632835800084
$ws.Range("J5").Value = 0.5080632835800084
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 3.118221666666667
$ws.Range("N5").Value = 9.354665000000001
$ws.Range("O5").Value = 0.01365525192599884
$ws.Range("P5").Value = 0.01365525192599884
$ws.Range("Q5").Value = 662.3340937800539
$ws.Range("R5").Value = 5961.006844020485
$ws.Range("S5").Value = 0.006937732131635206
$ws.Range("T5").Value = 0.006937732131635206

$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 161.9384456666667
$ws.Range("H6").Value = 485.815337
$ws.Range("I6").Value = 0.3873447295187379
$ws.Range("J6").Value = 0.3873447295187379
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 169.629438
$ws.Range("N6").Value = 508.888314
$ws.Range("O6").Value = 0.7428377317484701
$ws.Range("P6").Value = 0.7428377317484702
$ws.Range("Q6").Value = 27469.5275290302
$ws.Range("R6").Value = 247225.7477612718
$ws.Range("S6").Value = 0.287734280280424
$ws.Range("T6").Value = 0.287734280280424

$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 161.9384456666667
$ws.Range("H7").Value = 485.815337
$ws.Range("I7").Value = 0.3873447295187379
$ws.Range("J7").Value = 0.3873447295187379
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 0.9848756666666668
$ws.Range("N7").Value = 2.954627
$ws.Range("O7").Value = 0.004312947180081616
$ws.Range("P7").Value = 0.004312947180081616
$ws.Range("Q7").Value = 159.4892346349221
$ws.Range("R7").Value = 1435.403111714299
$ws.Range("S7").Value = 0.001670597358897317
$ws.Range("T7").Value = 0.001670597358897317

$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 161.9384456666667
$ws.Range("H8").Value = 485.815337
$ws.Range("I8").Value = 0.3873447295187379
$ws.Range("J8").Value = 0.3873447295187379
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 54.620752
$ws.Range("N8").Value = 163.862256
$ws.Range("O8").Value = 0.2391940691454494
$ws.Range("P8").Value = 0.2391940691454494
$ws.Range("Q8").Value = 8845.199680024476
$ws.Range("R8").Value = 79606.79712022028
$ws.Range("S8").Value = 0.09265056201563039
$ws.Range("T8").Value = 0.09265056201563039

$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 161.9384456666667
$ws.Range("H9").Value = 485.815337
$ws.Range("I9").Value = 0.3873447295187379
$ws.Range("J9").Value = 0.3873447295187379
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 3.118221666666667
$ws.Range("N9").Value = 9.354665000000001
$ws.Range("O9").Value = 0.01365525192599884
$ws.Range("P9").Value = 0.01365525192599884
$ws.Range("Q9").Value = 504.9599699441228
$ws.Range("R9").Value = 4544.639729497106
$ws.Range("S9").Value = 0.005289289863786246
$ws.Range("T9").Value = 0.005289289863786246

$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.338549
$ws.Range("H10").Value = 1.015647
$ws.Range("I10").Value = 0.0008097840527861261
$ws.Range("J10").Value = 0.0008097840527861261
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 169.629438
$ws.Range("N10").Value = 508.888314
$ws.Range("O10").Value = 0.7428377317484701
$ws.Range("P10").Value = 0.7428377317484702
$ws.Range("Q10").Value = 57.427876605462
$ws.Range("R10").Value = 516.850889449158
$ws.Range("S10").Value = 0.0006015381489777293
$ws.Range("T10").Value = 0.0006015381489777293

$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.338549
$ws.Range("H11").Value = 1.015647
$ws.Range("I11").Value = 0.0008097840527861261
$ws.Range("J11").Value = 0.0008097840527861261
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 0.9848756666666668
$ws.Range("N11").Value = 2.954627
$ws.Range("O11").Value = 0.004312947180081616
$ws.Range("P11").Value = 0.004312947180081616
$ws.Range("Q11").Value = 0.3334286720743334
$ws.Range("R11").Value = 3.000858048669
$ws.Range("S11").Value = 0.000003492555846938985
$ws.Range("T11").Value = 0.000003492555846938985

$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.338549
$ws.Range("H12").Value = 1.015647
$ws.Range("I12").Value = 0.0008097840527861261
$ws.Range("J12").Value = 0.0008097840527861261
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 54.620752
$ws.Range("N12").Value = 163.862256
$ws.Range("O12").Value = 0.2391940691454494
$ws.Range("P12").Value = 0.2391940691454494
$ws.Range("Q12").Value = 18.491800968848
$ws.Range("R12").Value = 166.426208719632
$ws.Range("S12").Value = 0.0001936955427150069
$ws.Range("T12").Value = 0.0001936955427150069

$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.338549
$ws.Range("H13").Value = 1.015647
$ws.Range("I13").Value = 0.0008097840527861261
$ws.Range("J13").Value = 0.0008097840527861261
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 3.118221666666667
$ws.Range("N13").Value = 9.354665000000001
$ws.Range("O13").Value = 0.01365525192599884
$ws.Range("P13").Value = 0.01365525192599884
$ws.Range("Q13").Value = 1.055670827028333
$ws.Range("R13").Value = 9.501037443255001
$ws.Range("S13").Value = 0.0000110578052464509
$ws.Range("T13").Value = 0.0000110578052464509

$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 43.38855633333333
$ws.Range("H14").Value = 130.165669
$ws.Range("I14").Value = 0.1037822028484675
$ws.Range("J14").Value = 0.1037822028484675
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 169.629438
$ws.Range("N14").Value = 508.888314
$ws.Range("O14").Value = 0.7428377317484701
$ws.Range("P14").Value = 0.7428377317484702
$ws.Range("Q14").Value = 7359.976426454674
$ws.Range("R14").Value = 66239.78783809207
$ws.Range("S14").Value = 0.07709333615981517
$ws.Range("T14").Value = 0.07709333615981519

$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 43.38855633333333
$ws.Range("H15").Value = 130.165669
$ws.Range("I15").Value = 0.1037822028484675
$ws.Range("J15").Value = 0.1037822028484675
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 0.9848756666666668
$ws.Range("N15").Value = 2.954627
$ws.Range("O15").Value = 0.004312947180081616
$ws.Range("P15").Value = 0.004312947180081616
$ws.Range("Q15").Value = 42.7323333444959
$ws.Range("R15").Value = 384.5910001004631
$ws.Range("S15").Value = 0.0004476071591179559
$ws.Range("T15").Value = 0.0004476071591179559

$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 43.38855633333333
$ws.Range("H16").Value = 130.165669
$ws.Range("I16").Value = 0.1037822028484675
$ws.Range("J16").Value = 0.1037822028484675
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 54.620752
$ws.Range("N16").Value = 163.862256
$ws.Range("O16").Value = 0.2391940691454494
$ws.Range("P16").Value = 0.2391940691454494
$ws.Range("Q16").Value = 2369.91557512103
$ws.Range("R16").Value = 21329.24017608926
$ws.Range("S16").Value = 0.02482408740420338
$ws.Range("T16").Value = 0.02482408740420338

$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 43.38855633333333
$ws.Range("H17").Value = 130.165669
$ws.Range("I17").Value = 0.1037822028484675
$ws.Range("J17").Value = 0.1037822028484675
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 3.118221666666667
$ws.Range("N17").Value = 9.354665000000001
$ws.Range("O17").Value = 0.01365525192599884
$ws.Range("P17").Value = 0.01365525192599884
$ws.Range("Q17").Value = 135.2951364439872
$ws.Range("R17").Value = 1217.656227995885
$ws.Range("S17").Value = 0.001417172125330938
$ws.Range("T17").Value = 0.001417172125330938
